$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 601
$ws.Range("I5").Value = 461
$ws.Range("J5").Value = 951
$ws.Range("K5").Value = 461
$ws.Range("L5").Value = 951
$ws.Range("M5").Value = -346
$ws.Range("N5").Value = -1181

# Row 62
$ws.Range("H62").Value = 2188.7778
$ws.Range("I62").Value = 1949.8334
$ws.Range("J62").Value = 2666.6667
$ws.Range("K62").Value = 1949.8334
$ws.Range("L62").Value = 2666.6667
$ws.Range("M62").Value = -1325.8334
$ws.Range("N62").Value = -3914.6667

# Row 65
$ws.Range("H65").Value = 2188.7778
$ws.Range("I65").Value = 1949.8334
$ws.Range("J65").Value = 2666.6667
$ws.Range("K65").Value = 9749.166999999999
$ws.Range("L65").Value = 13333.3335
$ws.Range("M65").Value = -6629.166999999999
$ws.Range("N65").Value = -19573.3335

# Row 95
$ws.Range("H95").Value = 50666.668
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 50666.668
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 50666.668
$ws.Range("N95").Value = -56158.668

# Row 127
$ws.Range("H127").Value = 813.5925999999999
$ws.Range("I127").Value = 349.8
$ws.Range("J127").Value = 1393.3334
$ws.Range("K127").Value = 1049.4
$ws.Range("L127").Value = 4180.0002
$ws.Range("M127").Value = 3910.6
$ws.Range("N127").Value = -14100.0002

# Row 138
$ws.Range("H138").Value = 1419.49
$ws.Range("I138").Value = 620.8378
$ws.Range("J138").Value = 1888.5397
$ws.Range("K138").Value = 1862.5134
$ws.Range("L138").Value = 5665.6191
$ws.Range("M138").Value = 3277.4866
$ws.Range("N138").Value = -15945.6191


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 679238.0600000001
$ws.Range("I32").Value = 796886.6
$ws.Range("J32").Value = 15364.071
$ws.Range("K32").Value = 796886.6
$ws.Range("L32").Value = 15364.071
$ws.Range("M32").Value = -796599.6
$ws.Range("N32").Value = -15938.071

# Row 132
$ws.Range("H132").Value = 3099.8838
$ws.Range("I132").Value = 2844.138
$ws.Range("J132").Value = 3629.6428
$ws.Range("K132").Value = 8532.414000000001
$ws.Range("L132").Value = 10888.9284
$ws.Range("M132").Value = -6002.414000000001
$ws.Range("N132").Value = -15948.9284


$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1259.9
$ws.Range("I58").Value = 953.46875
$ws.Range("J58").Value = 1804.6666
$ws.Range("K58").Value = 953.46875
$ws.Range("L58").Value = 1804.6666
$ws.Range("M58").Value = -750.46875
$ws.Range("N58").Value = -2210.6666

# Row 99
$ws.Range("H99").Value = 2072.56
$ws.Range("I99").Value = 1940
$ws.Range("J99").Value = 2078.0833
$ws.Range("K99").Value = 1940
$ws.Range("L99").Value = 2078.0833
$ws.Range("M99").Value = -442
$ws.Range("N99").Value = -5074.0833

# Row 105
$ws.Range("H105").Value = 1024.5294
$ws.Range("I105").Value = 1025.5
$ws.Range("J105").Value = 1022.2
$ws.Range("K105").Value = 1025.5
$ws.Range("L105").Value = 1022.2
$ws.Range("M105").Value = 721.5
$ws.Range("N105").Value = -4516.2

# Row 126
$ws.Range("H126").Value = 2072.56
$ws.Range("I126").Value = 1940
$ws.Range("J126").Value = 2078.0833
$ws.Range("K126").Value = 5820
$ws.Range("L126").Value = 6234.249899999999
$ws.Range("M126").Value = -3350
$ws.Range("N126").Value = -11174.2499

# Row 134
$ws.Range("H134").Value = 3684
$ws.Range("I134").Value = 4061.3333
$ws.Range("J134").Value = 2726.1538
$ws.Range("K134").Value = 12183.9999
$ws.Range("L134").Value = 8178.4614
$ws.Range("M134").Value = -9648.999899999999
$ws.Range("N134").Value = -13248.4614

# Row 136
$ws.Range("H136").Value = 1259.9
$ws.Range("I136").Value = 953.46875
$ws.Range("J136").Value = 1804.6666
$ws.Range("K136").Value = 2860.40625
$ws.Range("L136").Value = 5413.9998
$ws.Range("M136").Value = -310.40625
$ws.Range("N136").Value = -10513.9998


$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 5206
$ws.Range("I56").Value = 5206
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 5206
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -4676

# Row 64
$ws.Range("H64").Value = 1818.5714
$ws.Range("I64").Value = 863
$ws.Range("J64").Value = 2200.8
$ws.Range("K64").Value = 2589
$ws.Range("L64").Value = 6602.400000000001
$ws.Range("M64").Value = -2319
$ws.Range("N64").Value = -7142.400000000001

# Row 67
$ws.Range("H67").Value = 1818.5714
$ws.Range("I67").Value = 863
$ws.Range("J67").Value = 2200.8
$ws.Range("K67").Value = 2589
$ws.Range("L67").Value = 6602.400000000001
$ws.Range("M67").Value = -1653
$ws.Range("N67").Value = -8474.400000000001

# Row 75
$ws.Range("H75").Value = 3182.6
$ws.Range("I75").Value = 191.5
$ws.Range("J75").Value = 5176.6665
$ws.Range("K75").Value = 574.5
$ws.Range("L75").Value = 15529.9995
$ws.Range("M75").Value = 423.5
$ws.Range("N75").Value = -17525.9995

# Row 78
$ws.Range("H78").Value = 3182.6
$ws.Range("I78").Value = 191.5
$ws.Range("J78").Value = 5176.6665
$ws.Range("K78").Value = 1723.5
$ws.Range("L78").Value = 46589.9985
$ws.Range("M78").Value = 3268.5
$ws.Range("N78").Value = -56573.9985

# Row 80
$ws.Range("H80").Value = 3292.5
$ws.Range("I80").Value = 4202
$ws.Range("J80").Value = 3162.5715
$ws.Range("K80").Value = 12606
$ws.Range("L80").Value = 9487.7145
$ws.Range("M80").Value = -11670
$ws.Range("N80").Value = -11359.7145

# Row 81
$ws.Range("H81").Value = 8555.444
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 9374.875
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 28124.625
$ws.Range("M81").Value = -4877
$ws.Range("N81").Value = -30370.625

# Row 83
$ws.Range("H83").Value = 3292.5
$ws.Range("I83").Value = 4202
$ws.Range("J83").Value = 3162.5715
$ws.Range("K83").Value = 37818
$ws.Range("L83").Value = 28463.1435
$ws.Range("M83").Value = -33138
$ws.Range("N83").Value = -37823.1435

# Row 84
$ws.Range("H84").Value = 8555.444
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 9374.875
$ws.Range("K84").Value = 18000
$ws.Range("L84").Value = 84373.875
$ws.Range("M84").Value = -12384
$ws.Range("N84").Value = -95605.875

# Row 87
$ws.Range("H87").Value = 966.6667
$ws.Range("I87").Value = 966.6667
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 2900.0001
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -1652.0001
$ws.Range("N87").Value = ""

# Row 90
$ws.Range("H90").Value = 966.6667
$ws.Range("I90").Value = 966.6667
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 8700.0003
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -2460.0003
$ws.Range("N90").Value = ""

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""

# Row 95
$ws.Range("H95").Value = 9683.666999999999
$ws.Range("I95").Value = 1024
$ws.Range("J95").Value = 14013.5
$ws.Range("K95").Value = 3072
$ws.Range("L95").Value = 42040.5
$ws.Range("M95").Value = -1013
$ws.Range("N95").Value = -46158.5

# Row 103
$ws.Range("H103").Value = 128.6
$ws.Range("I103").Value = 128.6
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 385.8
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 493.2
$ws.Range("N103").Value = ""

# Row 106
$ws.Range("H106").Value = 8361.809999999999
$ws.Range("I106").Value = 3026
$ws.Range("J106").Value = 8628.6
$ws.Range("K106").Value = 9078
$ws.Range("L106").Value = 25885.8
$ws.Range("M106").Value = -8132
$ws.Range("N106").Value = -27777.8

# Row 107
$ws.Range("H107").Value = 21276912
$ws.Range("I107").Value = 292.7097
$ws.Range("J107").Value = 62500360
$ws.Range("K107").Value = 878.1291
$ws.Range("L107").Value = 187501080
$ws.Range("M107").Value = 1041.8709
$ws.Range("N107").Value = -187504920

# Row 113
$ws.Range("H113").Value = 479.03333
$ws.Range("I113").Value = 476
$ws.Range("J113").Value = 481.6875
$ws.Range("K113").Value = 1428
$ws.Range("L113").Value = 1445.0625
$ws.Range("M113").Value = 742
$ws.Range("N113").Value = -5785.0625

# Row 114
$ws.Range("H114").Value = 1363.375
$ws.Range("I114").Value = 300
$ws.Range("J114").Value = 1717.8334
$ws.Range("K114").Value = 900
$ws.Range("L114").Value = 5153.5002
$ws.Range("M114").Value = 2354
$ws.Range("N114").Value = -11661.5002

# Row 117
$ws.Range("H117").Value = 522.4
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 522.4
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 1567.2
$ws.Range("N117").Value = -8451.200000000001

# Row 131
$ws.Range("H131").Value = 2849.0476
$ws.Range("I131").Value = 465.9
$ws.Range("J131").Value = 3298.698
$ws.Range("K131").Value = 1397.7
$ws.Range("L131").Value = 9896.093999999999
$ws.Range("M131").Value = 3642.3
$ws.Range("N131").Value = -19976.094

# Row 137
$ws.Range("H137").Value = 6179294.5
$ws.Range("I137").Value = 20847130
$ws.Range("J137").Value = 3363.8948
$ws.Range("K137").Value = 62541390
$ws.Range("L137").Value = 10091.6844
$ws.Range("M137").Value = -62536290
$ws.Range("N137").Value = -20291.6844


$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1693.5657
$ws.Range("I132").Value = 1318.0844
$ws.Range("J132").Value = 3641.375
$ws.Range("K132").Value = 3954.2532
$ws.Range("L132").Value = 10924.125
$ws.Range("M132").Value = -1424.2532
$ws.Range("N132").Value = -15984.125


$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2359.1724
$ws.Range("I82").Value = 2331.7896
$ws.Range("J82").Value = 2411.2
$ws.Range("K82").Value = 2331.7896
$ws.Range("L82").Value = 2411.2
$ws.Range("M82").Value = -1970.7896
$ws.Range("N82").Value = -3133.2

# Row 85
$ws.Range("H85").Value = 2359.1724
$ws.Range("I85").Value = 2331.7896
$ws.Range("J85").Value = 2411.2
$ws.Range("K85").Value = 2331.7896
$ws.Range("L85").Value = 2411.2
$ws.Range("M85").Value = -1083.7896
$ws.Range("N85").Value = -4907.2

# Row 122
$ws.Range("H122").Value = 4471.1
$ws.Range("I122").Value = 3741.9167
$ws.Range("J122").Value = 4957.222
$ws.Range("K122").Value = 11225.7501
$ws.Range("L122").Value = 14871.666
$ws.Range("M122").Value = -8775.750100000001
$ws.Range("N122").Value = -19771.666


$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1339.7693
$ws.Range("I126").Value = 1359.3334
$ws.Range("J126").Value = 1105
$ws.Range("K126").Value = 4078.0002
$ws.Range("L126").Value = 3315
$ws.Range("M126").Value = -1608.0002
$ws.Range("N126").Value = -8255

# Row 132
$ws.Range("H132").Value = 4862726.5
$ws.Range("I132").Value = 1795.8182
$ws.Range("J132").Value = 10803864
$ws.Range("K132").Value = 5387.4546
$ws.Range("L132").Value = 32411592
$ws.Range("M132").Value = -2857.4546
$ws.Range("N132").Value = -32416652

# Row 136
$ws.Range("H136").Value = 1802.45
$ws.Range("I136").Value = 1724.1892
$ws.Range("J136").Value = 2025.1923
$ws.Range("K136").Value = 5172.5676
$ws.Range("L136").Value = 6075.5769
$ws.Range("M136").Value = -2622.5676
$ws.Range("N136").Value = -11175.5769

